$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") contains a date serial number that changed from
# 45182 (2023-09-13) to 45184 (2023-09-15) for every data row (rows 2-260).
$newDate = 45184

for ($row = 2; $row -le 260; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
